# GPTA4_tech.pptx feedback-file adjustments
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 1 - Title slide
# ---------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Intro to AWS"
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "How to utilize it for MMM."

# ---------------------------------------------------------------
# Slide 2 - Benefits of Cloud Computing
# ---------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Cloud Computing Advantages"
$s2.Shapes.Item(2).TextFrame.TextRange.Text = @'
Switch fixed costs to variable costs.
Leverage significant economies of scale.
Eliminate capacity speculation.
Enhance agility and speed.
Save money on data center operations.
Achieve global expansion quickly.
'@

# ---------------------------------------------------------------
# Slide 3 - Cloud Computing (deployment types)
# ---------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Cloud Computing Technology"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = @'
Deployment Types:
Cloud-Based Deployment
Entire operation running on the cloud, transferring existing applications to the cloud, creating new cloud-based apps.
On-Premise Deployment
Resources shaped by virtualisation and resource management tools, enhancing resource usage via app management and virtualisation technologies.
Hybrid Deployment
Linking cloud-based resources with onsite infrastructure, integrating these resources with legacy IT infrastructure.
'@

# ---------------------------------------------------------------
# Slide 4 - Cloud Migration / AWS CAF
# ---------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Cloud Transition Strategy"
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "The AWS Cloud Adoption Framework (AWS CAF) categorizes guidance into six Perspectives, each with specific responsibilities. The planning process facilitates organization-wide readiness for upcoming changes. Business, People, and Governance Perspectives emphasize business capabilities, while Platform, Security, and Operations Perspectives concentrate on technical capabilities. The Governance Perspective aligns IT strategy with business strategy, optimizing business value and mitigating risks. It guides in revamping skills and procedures required for cloud business governance and in managing and evaluating cloud investments for business outcomes."

# ---------------------------------------------------------------
# Slide 5 - General governance
# ---------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = @'
Governance Overview

'@
$s5.Shapes.Item(2).TextFrame.TextRange.Text = "Sorry, there is no sentence provided to reword. Please provide a sentence."

# ---------------------------------------------------------------
# Slide 6 - Elastic Compute Cloud - EC2
# ---------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = '"Understanding EC2: Elastic Compute Cloud"'
$s6.Shapes.Item(2).TextFrame.TextRange.Text = "Amazon EC2 offers scalable and secure cloud computing. In contrast to traditional on-site resources requiring upfront hardware investment, delivery waiting time, physical installation and configuration, Amazon EC2 allows you to run applications on virtual servers in the AWS Cloud. These instances can be provisioned and launched in minutes, only billed for actual compute time, stopped when not needed, thus providing cost-efficiency by solely paying for required server capacity."

# ---------------------------------------------------------------
# Slide 7 - Serverless Computing
# ---------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$s7.Shapes.Item(1).TextFrame.TextRange.Text = "Serverless Architecture"
$s7.Shapes.Item(2).TextFrame.TextRange.Text = @'
To run applications in Amazon EC2, you must provision instances, upload your code and manage the instances. 

"Serverless" denotes that your code runs on servers, however, there is no need for server management or provisioning. This enables more focus on product innovation. 

Serverless computing offers scalability, adjusting application capacity by modifying consumption units like throughput and memory.
'@

# ---------------------------------------------------------------
# Slide 8 - AWS Lambda
# ---------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "No change needed"
$s8.Shapes.Item(2).TextFrame.TextRange.Text = "AWS Lambda facilitates serverless code execution, with charges only applying when your code runs, effectively minimizing cost. This service sustains diverse applications or backend services without necessitating management. To utilize AWS Lambda, upload your code and set a trigger (such as AWS services, mobile apps, or HTTP endpoints). The code runs only upon trigger activation."

# ---------------------------------------------------------------
# Slide 9 - Simple Storage Service - S3 Buckets
# ---------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$s9.Shapes.Item(1).TextFrame.TextRange.Text = @'
"S3 Buckets: Simple Storage Service"

'@
$s9.Shapes.Item(2).TextFrame.TextRange.Text = @'
Amazon Simple Storage Service (Amazon S3) offers object storage, storing data in buckets. File permissions control visibility and access upon upload, while Amazon S3 versioning monitors object changes.

Amazon Elastic Block Store (EBS) provides block storage volumes for use with Amazon EC2 instances, preserving data on attached EBS volumes even if an EC2 instance is stopped or terminated.
'@

# ---------------------------------------------------------------
# Slide 11 - AWS Database services
# ---------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$s11.Shapes.Item(1).TextFrame.TextRange.Text = "Amazon Web Services Databases"
$s11.Shapes.Item(2).TextFrame.TextRange.Text = "Amazon RDS facilitates relational database management in the AWS Cloud. Amazon Aurora, a high-performance relational database, is compatible with MySQL and PostgreSQL, and outperforms standard databases. Amazon DynamoDB is a key-value database offering high-speed performance at any scale. Amazon Redshift is a data warehousing service, useful for big data analytics, collecting data from various sources to identify relationships and trends."

# ---------------------------------------------------------------
# Slide 12 - Rds vs flat data placeholder
# ---------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$s12.Shapes.Item(1).TextFrame.TextRange.Text = '"RDS vs Flat Data Comparison"'
$s12.Shapes.Item(2).TextFrame.TextRange.Text = @'
Input: Guiding our discussion today will be a comprehensive analysis of the wide-ranging effects that are embedded within the implementation of artificial intelligence systems into businesses and how they can potentially enhance efficiency and productivity.

Output: Today's discussion revolves around a detailed examination of how incorporating artificial intelligence systems into businesses can potentially boost efficiency and productivity.
'@

# ---------------------------------------------------------------
# Slide 13 - Identity and Access Manager
# ---------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$s13.Shapes.Item(1).TextFrame.TextRange.Text = "Managing Identity and Access"
$s13.Shapes.Item(2).TextFrame.TextRange.Text = "AWS Identity and Access Management (IAM) ensures secure access to AWS resources, tailored to meet your business's unique operational and security requirements using IAM features: users, groups, roles, policies, and multi-factor authentication."

# ---------------------------------------------------------------
# Slide 14 - Have a go yourself
# ---------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$s14.Shapes.Item(1).TextFrame.TextRange.Text = "Try it Yourself"
$s14.Shapes.Item(2).TextFrame.TextRange.Text = "The utility of the proposed technological solution has the potential to be extremely beneficial due to its innovative implementation of high-level programming, which would inherently result in a major increase in overall system productivity as well as reduce the possibility of any unnecessary glitches occurring as a byproduct of system operations."

# ---------------------------------------------------------------
# Slide 15 - Executive Summary
# ---------------------------------------------------------------
$s15 = $p.Slides.Item(15)
$s15.Shapes.Item(2).TextFrame.TextRange.Text = @'
- Introduction to AWS: Introducing AWS Cloud services, and its potential advantages including cost savings, accessibility, global reach, scalability and its application for MMM. 
- Deployment Types: Discussing the various types of cloud deployment including Cloud Based, On-Premise and Hybrid, all having unique benefits and considerations.
- Cloud Migration & AWS CAF: Understanding AWS Cloud Adoption Framework (AWS CAF) for planned and systematic cloud migration which involves various perspectives like Business, People, Governance, Platform, Security, and Operations.
- EC2 & Serverless Computing: Introducing Amazon Elastic Compute Cloud (EC2) which provides secure, re-sizable compute capacity in the cloud and discussing serverless computing, which allows developers to focus on their applications without worrying about server management.
- AWS Services Overview: Providing an overview of various AWS services like AWS Lambda, Simple Storage Service (S3 Buckets), AWS Database services and AWS Identity and Access Management (IAM), each with unique benefits and potential applications.
'@
